$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Value) {
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "91.432.41"
Set-TextValue $ws.Range("E2") "  +4.14%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.185.15"
Set-TextValue $ws.Range("E3") "  +0.64%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "217.22"
Set-TextValue $ws.Range("E5") "  +4.64%  "

# Row 6
Set-TextValue $ws.Range("D6") "634.29"
Set-TextValue $ws.Range("E6") "  +4.26%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.392"
Set-TextValue $ws.Range("E7") "  +2.37%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.723"
Set-TextValue $ws.Range("E8") "  +7.79%  "

# Row 9
Set-TextValue $ws.Range("D9") "1.00"
Set-TextValue $ws.Range("E9") "  +0.07%  "

# Row 10
Set-TextValue $ws.Range("D10") "3.184.26"
Set-TextValue $ws.Range("E10") "  +0.78%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.567"
Set-TextValue $ws.Range("E11") "  +6.30%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.181"
Set-TextValue $ws.Range("E12") "  +2.87%  "

# Row 13
Set-TextValue $ws.Range("E13") "  +5.39%  "

# Row 14
Set-TextValue $ws.Range("D14") "91.013.27"
Set-TextValue $ws.Range("E14") "  +3.91%  "

# Row 15
Set-TextValue $ws.Range("D15") "5.33"
Set-TextValue $ws.Range("E15") "  +0.96%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.771.34"
Set-TextValue $ws.Range("E16") "  +0.81%  "

# Row 17
Set-TextValue $ws.Range("D17") "32.57"
Set-TextValue $ws.Range("E17") "  +1.36%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.191.40"
Set-TextValue $ws.Range("E18") "  +0.53%  "

# Row 19
Set-TextValue $ws.Range("D19") "3.30"
Set-TextValue $ws.Range("E19") "  +3.38%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0000212"
Set-TextValue $ws.Range("E20") "  +60.60%  "

# Row 21
Set-TextValue $ws.Range("D21") "436.24"
Set-TextValue $ws.Range("E21") "  +6.08%  "

# Row 22
Set-TextValue $ws.Range("D22") "13.28"
Set-TextValue $ws.Range("E22") "  -1.15%  "

# Row 23
Set-TextValue $ws.Range("D23") "8.49"
Set-TextValue $ws.Range("E23") "  +0.64%  "

# Row 24
Set-TextValue $ws.Range("D24") "4.99"
Set-TextValue $ws.Range("E24") "  -1.36%  "

# Row 25
Set-TextValue $ws.Range("D25") "5.21"
Set-TextValue $ws.Range("E25") "  -0.23%  "

# Row 26
Set-TextValue $ws.Range("D26") "11.64"
Set-TextValue $ws.Range("E26") "  -4.51%  "

# Row 27
Set-TextValue $ws.Range("D27") "80.35"
Set-TextValue $ws.Range("E27") "  +9.56%  "

# Row 28
Set-TextValue $ws.Range("D28") "3.349.93"
Set-TextValue $ws.Range("E28") "  +0.46%  "

# Row 29
Set-TextValue $ws.Range("E29") "  +0.11%  "

# Row 30
Set-TextValue $ws.Range("B30") "Cronos"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D30") "0.161"
Set-TextValue $ws.Range("E30") "  -0.52%  "

# Row 31
Set-TextValue $ws.Range("B31") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D31") "1.00"
Set-TextValue $ws.Range("E31") "  -0.07%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.01"
Set-TextValue $ws.Range("E32") "  +32.76%  "

# Row 33
Set-TextValue $ws.Range("E33") "  +1.36%  "

# Row 34
Set-TextValue $ws.Range("D34") "517.71"
Set-TextValue $ws.Range("E34") "  -5.24%  "

# Row 35
Set-TextValue $ws.Range("D35") "7.00"
Set-TextValue $ws.Range("E35") "  +1.15%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.89"
Set-TextValue $ws.Range("E36") "  +1.93%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -1.99%  "

# Row 38
Set-TextValue $ws.Range("D38") "22.34"
Set-TextValue $ws.Range("E38") "  +2.27%  "

# Row 39
Set-TextValue $ws.Range("D39") "22.43"
Set-TextValue $ws.Range("E39") "  +2.86%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +0.27%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.125"
Set-TextValue $ws.Range("E41") "  -3.00%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -0.06%  "

# Row 43
Set-TextValue $ws.Range("E43") "  +0.64%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.369"
Set-TextValue $ws.Range("E44") "  -0.82%  "

# Row 45
Set-TextValue $ws.Range("D45") "146.98"
Set-TextValue $ws.Range("E45") "  -2.23%  "

# Row 46
Set-TextValue $ws.Range("D46") "44.16"
Set-TextValue $ws.Range("E46") "  +2.06%  "

# Row 47
Set-TextValue $ws.Range("D47") "170.62"
Set-TextValue $ws.Range("E47") "  -1.66%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +1.37%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.742"
Set-TextValue $ws.Range("E49") "  +7.85%  "

# Row 50
Set-TextValue $ws.Range("D50") "24.62"
Set-TextValue $ws.Range("E50") "  +3.26%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.609"
Set-TextValue $ws.Range("E51") "  +3.27%  "
